# Apply updated cryptocurrency price/volume data to Sheet1 (rows 2-51).
# D column = Price (text, may look numeric -> must stay literal text)
# E column = Volume(1h) percent change (always text, "  +x.xx%  ")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Replace("68.320.84", "68.149.93") | Out-Null
$ws.Range("E2").Replace("  +8.63%  ", "  +8.26%  ") | Out-Null

$ws.Range("D3").Replace("3.631.72", "3.624.98") | Out-Null
$ws.Range("E3").Replace("  +4.54%  ", "  +4.21%  ") | Out-Null

$ws.Range("E4").Replace("  +0.38%  ", "  +0.11%  ") | Out-Null

$ws.Range("Z1").Value = "'420.08"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Replace("  +1.25%  ", "  +1.22%  ") | Out-Null

$ws.Range("Z1").Value = "'131.43"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Replace("  +0.07%  ", "  +0.59%  ") | Out-Null

$ws.Range("Z1").Value = "'0.649"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Replace("  +3.46%  ", "  +3.34%  ") | Out-Null

$ws.Range("D8").Replace("3.624.28", "3.617.13") | Out-Null
$ws.Range("E8").Replace("  +4.56%  ", "  +4.22%  ") | Out-Null

$ws.Range("Z1").Value = "'1.00"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Replace("  -0.01%  ", "  +0.02%  ") | Out-Null

$ws.Range("Z1").Value = "'0.773"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Replace("  +6.48%  ", "  +6.08%  ") | Out-Null

$ws.Range("E11").Replace("  +20.30%  ", "  +18.37%  ") | Out-Null

$ws.Range("Z1").Value = "'0.0000356"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Replace("  +57.39%  ", "  +58.86%  ") | Out-Null

$ws.Range("Z1").Value = "'42.71"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Replace("  +0.34%  ", "  +0.04%  ") | Out-Null

$ws.Range("Z1").Value = "'9.93"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Replace("  +2.43%  ", "  +0.96%  ") | Out-Null

$ws.Range("D15").Replace("4.197.98", "4.200.02") | Out-Null
$ws.Range("E15").Replace("  +4.28%  ", "  +4.24%  ") | Out-Null

$ws.Range("Z1").Value = "'20.44"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Replace("  -0.26%  ", "  -0.31%  ") | Out-Null

$ws.Range("D18").Replace("3.608.69", "3.623.03") | Out-Null
$ws.Range("E18").Replace("  +3.71%  ", "  +3.70%  ") | Out-Null

$ws.Range("E19").Replace("  +4.11%  ", "  +4.06%  ") | Out-Null

$ws.Range("D20").Replace("68.219.45", "68.099.67") | Out-Null
$ws.Range("E20").Replace("  +8.59%  ", "  +8.29%  ") | Out-Null

$ws.Range("Z1").Value = "'12.45"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Replace("  -1.13%  ", "  -1.83%  ") | Out-Null

$ws.Range("Z1").Value = "'467.08"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Replace("  -1.30%  ", "  -0.88%  ") | Out-Null

$ws.Range("Z1").Value = "'88.96"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Replace("  -1.73%  ", "  -2.02%  ") | Out-Null

$ws.Range("E24").Replace("  -5.31%  ", "  -5.66%  ") | Out-Null

$ws.Range("Z1").Value = "'13.33"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Replace("  +1.00%  ", "  +1.05%  ") | Out-Null

$ws.Range("E26").Replace("  +0.98%  ", "  +0.85%  ") | Out-Null

$ws.Range("Z1").Value = "'10.10"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Replace("  -3.85%  ", "  -3.90%  ") | Out-Null

$ws.Range("Z1").Value = "'35.84"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Replace("  +7.61%  ", "  +6.93%  ") | Out-Null

$ws.Range("Z1").Value = "'4.88"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Replace("  +1.99%  ", "  +1.68%  ") | Out-Null

$ws.Range("E30").Replace("  +3.72%  ", "  +3.52%  ") | Out-Null

$ws.Range("Z1").Value = "'12.38"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Replace("  +1.86%  ", "  +1.72%  ") | Out-Null

$ws.Range("Z1").Value = "'0.117"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Replace("  +4.09%  ", "  +3.86%  ") | Out-Null

$ws.Range("Z1").Value = "'7.39"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Replace("  -2.42%  ", "  -2.23%  ") | Out-Null

$ws.Range("E34").Replace("  -3.94%  ", "  -3.78%  ") | Out-Null

$ws.Range("Z1").Value = "'40.69"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Replace("  -0.89%  ", "  -0.71%  ") | Out-Null

$ws.Range("Z1").Value = "'0.998"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Replace("  +0.07%  ", "  -0.14%  ") | Out-Null

$ws.Range("Z1").Value = "'56.80"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Replace("  -2.52%  ", "  -2.19%  ") | Out-Null

$ws.Range("Z1").Value = "'0.0495"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Replace("  +1.38%  ", "  +1.16%  ") | Out-Null

$ws.Range("D39").Replace("0.0₃0713", "0.0₃0708") | Out-Null
$ws.Range("E39").Replace("  +22.47%  ", "  +21.46%  ") | Out-Null

$ws.Range("E40").Replace("  +7.75%  ", "  +7.50%  ") | Out-Null

$ws.Range("E41").Replace("  +0.00%  ", "  +0.06%  ") | Out-Null

$ws.Range("Z1").Value = "'3.04"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Replace("  -0.40%  ", "  -0.34%  ") | Out-Null

$ws.Range("E43").Replace("  -3.24%  ", "  -2.79%  ") | Out-Null

$ws.Range("Z1").Value = "'148.18"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Replace("  -1.14%  ", "  -1.41%  ") | Out-Null

$ws.Range("Z1").Value = "'3.27"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Replace("  -1.80%  ", "  -1.85%  ") | Out-Null

$ws.Range("Z1").Value = "'4.32"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Replace("  -2.57%  ", "  -2.69%  ") | Out-Null

$ws.Range("E47").Replace("  -4.10%  ", "  -4.14%  ") | Out-Null

$ws.Range("E48").Replace("  -3.71%  ", "  -3.80%  ") | Out-Null

$ws.Range("Z1").Value = "'2.33"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Replace("  -2.67%  ", "  -2.90%  ") | Out-Null

$ws.Range("Z1").Value = "'2.72"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Replace("  +17.09%  ", "  +16.70%  ") | Out-Null

$ws.Range("Z1").Value = "'15.63"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Replace("  -4.98%  ", "  -4.97%  ") | Out-Null

$ws.Range("Z1").Clear() | Out-Null
